$wb = $excel.ActiveWorkbook

# --- Update the daily conversion note on "Hoja1" ---
$wsNota = $wb.Worksheets.Item("Hoja1")
$nl = [char]10
$lineas = @(
    "Conversión del día 💰",
    "✅ Dólar paralelo: 68",
    "",
    "Binance",
    "✅ 1000 Bs = 4.16 = 16124.13 pesos",
    "✅ 16124.13 pesos = 4.15 = 966.29 Bs",
    "",
    "Promedio competencia",
    "✅ Tasa pesos: 20",
    "✅ Tasa Bs: 20",
    "✅ % Ganancia: 20%"
)
$texto = [string]::Join($nl, $lineas)
$wsNota.Range("A1").Value = $texto

# --- Update the rate figures on "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 240.199
$wsTasas.Range("O10").Value = 3873
$wsTasas.Range("N12").Value = 3888
$wsTasas.Range("O12").Value = 233
